# A new review/avaliação row was recorded and inserted at the top of the
# existing data list (row 15, right after the first 13 unchanged rows),
# pushing every subsequent row down by one. Replicate that with a native
# row insert on the worksheet, then populate the new row's cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15:33 down to 16:34, leaving a blank row 15 behind.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row with the new review's data.
$ws.Range("A15").Value = 5
$ws.Range("C15").Value = 46010.46702450232
$ws.Range("D15").Value = "NmNjY2Q0MzYtNWI4ZS00ODk3LTgyZDItNTFkNWMxNWFjYzA5OjU3MDE2"
